$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.359.43"
$ws.Range("E2").Value = "  +1.16%  "

$ws.Range("D3").Value = "3.573.29"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").Value = "  +11.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.408"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.06"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.74%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").Value = "3.570.70"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.52"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.25%  "

$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").Value = "4.239.13"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").Value = "96.171.24"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("E17").Value = "  +1.95%  "

$ws.Range("D18").Value = "3.562.04"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.80"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.520"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.89%  "

$ws.Range("E23").Value = "  -5.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "501.85"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.60%  "

$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("D29").Value = "3.764.70"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.152"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.03%  "

$ws.Range("E31").Value = "  -5.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.183"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.35"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "622.44"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.79"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.565"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.36%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.903"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0425"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.56%  "

$ws.Range("E48").Value = "  +1.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.39"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("E51").Value = "  +2.64%  "
